# Update gh-pages to output generated at 456a3b4
# This script updates the "F" column (count) values on the "展览", "演出"
# and "全部类型" worksheets to reflect newly generated data.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 144
$ws1.Range("F4").Value  = 2118
$ws1.Range("F5").Value  = 380
$ws1.Range("F6").Value  = 665
$ws1.Range("F8").Value  = 2100
$ws1.Range("F9").Value  = 4
$ws1.Range("F10").Value = 10906
$ws1.Range("F11").Value = 182
$ws1.Range("F13").Value = 293
$ws1.Range("F14").Value = 208
$ws1.Range("F15").Value = 10697
$ws1.Range("F17").Value = 1124
$ws1.Range("F18").Value = 4
$ws1.Range("F19").Value = 745
$ws1.Range("F20").Value = 5351
$ws1.Range("F21").Value = 77
$ws1.Range("F22").Value = 3385

# ---- Sheet "演出" ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 26

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 144
$ws4.Range("F4").Value  = 2118
$ws4.Range("F5").Value  = 380
$ws4.Range("F6").Value  = 665
$ws4.Range("F7").Value  = 26
$ws4.Range("F9").Value  = 2100
$ws4.Range("F11").Value = 4
$ws4.Range("F13").Value = 10906
$ws4.Range("F14").Value = 182
$ws4.Range("F16").Value = 293
$ws4.Range("F17").Value = 208
$ws4.Range("F18").Value = 10697
$ws4.Range("F20").Value = 1124
$ws4.Range("F21").Value = 4
$ws4.Range("F22").Value = 745
$ws4.Range("F23").Value = 5351
$ws4.Range("F24").Value = 77
$ws4.Range("F25").Value = 3385
